$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.844.08"
$ws.Range("E2").Value = "  -5.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.60"
$ws.Range("E3").Value = "  -4.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "280.94"
$ws.Range("E5").Value = "  -8.33%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5050"
$ws.Range("E7").Value = "  -5.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3512"
$ws.Range("E8").Value = "  -7.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.04"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06628"
$ws.Range("E10").Value = "  -9.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.05"
$ws.Range("E11").Value = "  -9.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8514"
$ws.Range("E12").Value = "  -5.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07845"
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.807.79"
$ws.Range("E14").Value = "  +66.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.028"
$ws.Range("E15").Value = "  -5.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.33"
$ws.Range("E16").Value = "  -9.00%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.01"
$ws.Range("E18").Value = "  -5.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008102"
$ws.Range("E19").Value = "  -6.35%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.934.77"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.750"
$ws.Range("E22").Value = "  -5.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.10"
$ws.Range("E23").Value = "  -6.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.114"
$ws.Range("E24").Value = "  -6.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.67"
$ws.Range("E25").Value = "  -5.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.148"
$ws.Range("E26").Value = "  -5.95%  "
$ws.Range("E27").Value = "  -3.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.88"
$ws.Range("E28").Value = "  -7.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.58"
$ws.Range("E29").Value = "  -6.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.281"
$ws.Range("E30").Value = "  -11.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.194"
$ws.Range("E31").Value = "  -12.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08751"
$ws.Range("E32").Value = "  -5.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04792"
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7350"
$ws.Range("E34").Value = "  -11.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  -7.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.841"
$ws.Range("E36").Value = "  -5.11%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.090"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.447"
$ws.Range("E39").Value = "  -8.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5339"
$ws.Range("E40").Value = "  -7.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01847"
$ws.Range("E41").Value = "  -7.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9824"
$ws.Range("E42").Value = "  -8.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.73"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.187"
$ws.Range("E44").Value = "  -6.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.183"
$ws.Range("E45").Value = "  -12.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4678"
$ws.Range("E46").Value = "  -5.42%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1373"
$ws.Range("E48").Value = "  -9.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.243"
$ws.Range("E49").Value = "  -8.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.68"
$ws.Range("E50").Value = "  -6.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05904"
$ws.Range("E51").Value = "  -3.98%  "
